# Append new scrape run (2025-10-29 06:29 JST): refreshes the top-ranked
# rows and drops the ones that fell out of the top 9, per the commit
# message "Append: 2025-10-29 06:29 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-10-29 06:29:55"

# New data for rows 2-10 (row 1 is the header, left untouched).
# Columns: A=取得日時 B=タイトル C=カテゴリ D=価格 E=締切 F=URL G=優先度スコア H=スキル概要
$data = @(
    @{ B = "【業務効率化】生成AIを活用したシステム開発依頼"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422740"; G = 445; H = "🔥AI,Ai ◆効率化,開発" },
    @{ B = "【急募】業務効率化・生成AI実装のAIエンジニアパートナー募集"; D = "500,000 円 ~ 1,000,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422760"; G = 385; H = "🔥AI,Ai ◆効率化" },
    @{ B = "【機密性の高いノウハウを含む】サーバーレスAI分析システム構築(MVP開発と拡張性確保)"; D = "100,000 円 ~ 200,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422386"; G = 383; H = "🔥AI,Ai ◆開発" },
    @{ B = "【自動化】EAを証券口座・VPSに接続する作業の効率化依頼"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422331"; G = 145; H = "◆効率化,自動化" },
    @{ B = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"; D = "5,000 円 ~ 10,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5251319"; G = 135; H = "◆ツール,スクレイピング ◇サイト" },
    @{ B = "【動画解析】動作比較アルゴリズム開発者を募集します"; D = "50,000 円 ~ 100,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422314"; G = 68; H = "◆開発" },
    @{ B = "アマゾンの返品レポートより返品理由のポップアップ文字までダウンロードしてエクセルに書き出すツール"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422652"; G = 65; H = "◆ツール" },
    @{ B = "【オンライン講師募集】HTML・CSSの基礎を実践的に教えていただける方"; D = "50,000 円 ~ 100,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5422660"; G = 18; H = "" },
    @{ B = "【急募】YouTubeの音楽配信構築の依頼です"; D = "20,000 円 ~ 50,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5420233"; G = 13; H = "" }
)

# Drop the rows that fell out of the list (previously rows 11-21).
# Delete bottom-up so earlier row indices stay valid during the loop.
for ($r = 21; $r -ge 11; $r--) {
    $ws.Rows.Item($r).Delete()
}

# All existing hyperlinks point at stale URLs/relationship ids; clear the
# whole collection (any sub-range's .Hyperlinks.Delete() clears the sheet)
# and rebuild only the ones we still need, below.
$ws.Range("A1").Hyperlinks.Delete()

# Overwrite rows 2-10 with the refreshed data.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = "システム開発"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = "期限情報なし"
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G

    if ($row.H -ne "") {
        $ws.Cells.Item($r, 8).Value = $row.H
    } else {
        $ws.Cells.Item($r, 8).Value = $null
    }

    $linkCell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($linkCell, $row.F)
    $linkCell.Style = "Hyperlink"
}

# Column width tweaks from the diff (stored OOXML width = ColumnWidth + 5/6,
# so back the requested width off by 5/6 to land on the exact integer).
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
